# Weekly update: insert the latest "Naranja" price observation as a new
# row 84 (most-recent date first), pushing the existing rows 84-93 down
# to rows 85-94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 84, shifting rows 84:93 -> 85:94.
$ws.Rows.Item(84).Insert()

# Populate the new row 84 with this week's record.
$ws.Cells.Item(84, 1).Value  = 1
$ws.Cells.Item(84, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(84, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(84, 4).Value  = 44769
$ws.Cells.Item(84, 5).Value  = 15
$ws.Cells.Item(84, 6).Value  = "Fruta"
$ws.Cells.Item(84, 7).Value  = 100102
$ws.Cells.Item(84, 8).Value  = "Cítricos"
$ws.Cells.Item(84, 9).Value  = 100102005
$ws.Cells.Item(84, 10).Value = "Naranja"
$ws.Cells.Item(84, 11).Value = "New Hall"
$ws.Cells.Item(84, 12).Value = "Segunda"
$ws.Cells.Item(84, 13).Value = 270
$ws.Cells.Item(84, 14).Value = 650
$ws.Cells.Item(84, 15).Value = 700
$ws.Cells.Item(84, 16).Value = 675
$ws.Cells.Item(84, 17).Value = "`$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(84, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(84, 19).Value = 675
$ws.Cells.Item(84, 20).Value = 1
